$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($Cell, $Text)
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.Style = "Normal"
}

Set-TextCell $ws.Range("D2") "29.246.75"
Set-TextCell $ws.Range("E2") "  +1.50%  "

Set-TextCell $ws.Range("D3") "1.915.11"
Set-TextCell $ws.Range("E3") "  +1.75%  "

Set-TextCell $ws.Range("D4") "1.003"
Set-TextCell $ws.Range("E4") "  -0.24%  "

Set-TextCell $ws.Range("D5") "328.49"
Set-TextCell $ws.Range("E5") "  +1.60%  "

Set-TextCell $ws.Range("D6") "1.003"
Set-TextCell $ws.Range("E6") "  -0.01%  "

Set-TextCell $ws.Range("D7") "0.4665"
Set-TextCell $ws.Range("E7") "  -0.08%  "

Set-TextCell $ws.Range("D8") "0.3955"
Set-TextCell $ws.Range("E8") "  +0.50%  "

Set-TextCell $ws.Range("D9") "47.14"
Set-TextCell $ws.Range("E9") "  +1.31%  "

Set-TextCell $ws.Range("D10") "0.08044"
Set-TextCell $ws.Range("E10") "  +1.52%  "

Set-TextCell $ws.Range("D11") "1.019"
Set-TextCell $ws.Range("E11") "  +3.73%  "

Set-TextCell $ws.Range("D12") "22.29"
Set-TextCell $ws.Range("E12") "  -0.24%  "

Set-TextCell $ws.Range("D13") "1.901.93"
Set-TextCell $ws.Range("E13") "  +0.62%  "

Set-TextCell $ws.Range("D14") "7.180"
Set-TextCell $ws.Range("E14") "  +2.36%  "

Set-TextCell $ws.Range("D15") "5.807"
Set-TextCell $ws.Range("E15") "  +1.02%  "

Set-TextCell $ws.Range("D16") "0.06985"
Set-TextCell $ws.Range("E16") "  +0.15%  "

Set-TextCell $ws.Range("D17") "89.89"
Set-TextCell $ws.Range("E17") "  +1.22%  "

Set-TextCell $ws.Range("E18") "  -0.09%  "

Set-TextCell $ws.Range("D19") "0.00001020"
Set-TextCell $ws.Range("E19") "  +0.88%  "

Set-TextCell $ws.Range("D20") "17.39"
Set-TextCell $ws.Range("E20") "  +2.37%  "

Set-TextCell $ws.Range("E21") "  +0.18%  "

Set-TextCell $ws.Range("D22") "29.206.79"
Set-TextCell $ws.Range("E22") "  +1.32%  "

Set-TextCell $ws.Range("D23") "5.393"
Set-TextCell $ws.Range("E23") "  +0.74%  "

Set-TextCell $ws.Range("D24") "11.17"
Set-TextCell $ws.Range("E24") "  +0.65%  "

Set-TextCell $ws.Range("D25") "2.185.22"
Set-TextCell $ws.Range("E25") "  +3.04%  "

Set-TextCell $ws.Range("D26") "2.069"
Set-TextCell $ws.Range("E26") "  -2.24%  "

Set-TextCell $ws.Range("D27") "155.60"
Set-TextCell $ws.Range("E27") "  +1.24%  "

Set-TextCell $ws.Range("D28") "19.77"
Set-TextCell $ws.Range("E28") "  +1.86%  "

Set-TextCell $ws.Range("D29") "5.910"
Set-TextCell $ws.Range("E29") "  +2.37%  "

Set-TextCell $ws.Range("D30") "2.026"
Set-TextCell $ws.Range("E30") "  +1.16%  "

Set-TextCell $ws.Range("D31") "121.05"
Set-TextCell $ws.Range("E31") "  +0.77%  "

Set-TextCell $ws.Range("D32") "0.09399"
Set-TextCell $ws.Range("E32") "  +0.14%  "

Set-TextCell $ws.Range("D33") "0.9438"
Set-TextCell $ws.Range("E33") "  +0.38%  "

Set-TextCell $ws.Range("D34") "5.380"
Set-TextCell $ws.Range("E34") "  +1.17%  "

Set-TextCell $ws.Range("D35") "1.365"
Set-TextCell $ws.Range("E35") "  +0.24%  "

Set-TextCell $ws.Range("E36") "  -2.47%  "

Set-TextCell $ws.Range("D37") "0.05878"
Set-TextCell $ws.Range("E37") "  -0.72%  "

Set-TextCell $ws.Range("B38") "FraxShare"
Set-TextCell $ws.Range("C38") "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextCell $ws.Range("D38") "8.179"
Set-TextCell $ws.Range("E38") "  +3.48%  "

Set-TextCell $ws.Range("B39") "TrustWalletToken"
Set-TextCell $ws.Range("C39") "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextCell $ws.Range("D39") "1.181"
Set-TextCell $ws.Range("E39") "  +1.79%  "

Set-TextCell $ws.Range("D40") "0.02109"
Set-TextCell $ws.Range("E40") "  -0.85%  "

Set-TextCell $ws.Range("D41") "0.5849"
Set-TextCell $ws.Range("E41") "  +1.98%  "

Set-TextCell $ws.Range("E42") "  +0.04%  "

Set-TextCell $ws.Range("D43") "0.1824"

Set-TextCell $ws.Range("D44") "10.13"
Set-TextCell $ws.Range("E44") "  +1.33%  "

Set-TextCell $ws.Range("D45") "2.327"
Set-TextCell $ws.Range("E45") "  +9.40%  "

Set-TextCell $ws.Range("D46") "0.5483"
Set-TextCell $ws.Range("E46") "  +2.42%  "

Set-TextCell $ws.Range("D47") "11.98"
Set-TextCell $ws.Range("E47") "  +1.06%  "

Set-TextCell $ws.Range("D48") "0.07236"
Set-TextCell $ws.Range("E48") "  -0.91%  "

Set-TextCell $ws.Range("D49") "1.889"
Set-TextCell $ws.Range("E49") "  +2.13%  "

Set-TextCell $ws.Range("D50") "1.131"
Set-TextCell $ws.Range("E50") "  -3.81%  "

Set-TextCell $ws.Range("D51") "113.61"
Set-TextCell $ws.Range("E51") "  -0.66%  "
